$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting the existing row 3 (and below) down to row 4
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new match data
$ws.Range("A3").Value = "zNz4qPE7"
$ws.Range("B3").Value = "24/11/2024"
$ws.Range("C3").Value = "07:15"
$ws.Range("D3").Value = "BULGARIA - PARVA LIGA"
$ws.Range("E3").Value = "Beroe"
$ws.Range("F3").Value = "Septemvri Sofia"
$ws.Range("G3").Value = 1.67
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 2.38
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 6.5
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("S3").Value = 1.53
$ws.Range("T3").Value = 2.38
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
$ws.Range("W3").Value = 5
$ws.Range("X3").Value = 6.5
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 12
$ws.Range("AA3").Value = 17
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 6.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 101
$ws.Range("AH3").Value = 12
$ws.Range("AI3").Value = 29
$ws.Range("AJ3").Value = 21
$ws.Range("AK3").Value = 67
$ws.Range("AL3").Value = 51
$ws.Range("AM3").Value = 67
$ws.Range("AN3").Value = 3.4
$ws.Range("AO3").Value = 9
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 34
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 251
$ws.Range("AT3").Value = 2.38
$ws.Range("AU3").Value = 10
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 7
$ws.Range("AX3").Value = 34
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 151
$ws.Range("BA3").Value = 201
$ws.Range("BB3").Value = 51
$ws.Range("BC3").Value = 51
$ws.Range("BD3").Value = 51
